$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows right after the header block ends and the first
# existing data block of rows (77:80), shifting all existing data rows
# (formerly 77-176) down to 81-180.
$ws.Rows("77:80").Insert()

# Constant column values shared by every data row in this sheet.
$mercado = "Agr" + [char]0x00ED + "cola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$tipo = "Fruta"
$producto = "Tropicales y subtropicales"
$categoria = "Pi" + [char]0x00F1 + "a"
$variedad = "Caramelo"
$origen = "Ecuador"

# New weekly data block (date serial 44665) inserted at rows 77-80.
$rows = @(
    @{ Row = 77; Calidad = "Especial"; Volumen = 200; PMin = 18000; PMax = 19000; PProm = 18500; PKg = 1850; Kg = 10 },
    @{ Row = 78; Calidad = "Primera";  Volumen = 270; PMin = 18000; PMax = 19000; PProm = 18500; PKg = 1542; Kg = 12 },
    @{ Row = 79; Calidad = "Segunda";  Volumen = 250; PMin = 18000; PMax = 19000; PProm = 18500; PKg = 1321; Kg = 14 },
    @{ Row = 80; Calidad = "Tercera";  Volumen = 270; PMin = 18000; PMax = 19000; PProm = 18500; PKg = 1156; Kg = 16 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = 44665
    $ws.Cells.Item($row, 5).Value = 15
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = 100108
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = 100108005
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = "$/caja " + $r.Kg + " unidades"
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.Kg
}
